# Hortaliza, Vega Monumental Concepción - Repollo
# New weekly data point: insert two new rows (Primera / Segunda quality) at the
# top of the data block for this market/variety, pushing the rest of the
# historical rows down by two. The dataset was already sorted with the newest
# report rows on top of the block, so this mirrors that pattern. The two
# rows that previously fell off the bottom simply slide down with everything
# else (no separate append is required).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing data (rows 264-337) down by 2 rows, inserting two blank
# rows at 264-265.
$ws.Rows("264:265").Insert()

# --- New row 264: Repollo, Crespo record, Primera -------------------------
$ws.Range("A264").Value = 11
$ws.Range("B264").Value = 'Vega Monumental Concepción'
$ws.Range("C264").Value = 'Bíobío'
$ws.Range("D264").Value = 44736
$ws.Range("E264").Value = 8
$ws.Range("F264").Value = 100112006
$ws.Range("G264").Value = 'Repollo'
$ws.Range("H264").Value = 'Crespo record'
$ws.Range("I264").Value = 'Primera'
$ws.Range("J264").Value = 2000
$ws.Range("K264").Value = 800
$ws.Range("L264").Value = 900
$ws.Range("M264").Value = 850
$ws.Range("N264").Value = '$/unidad'
$ws.Range("O264").Value = 'Región Metropolitana'
$ws.Range("P264").Value = 850
$ws.Range("Q264").Value = 1
$ws.Range("R264").Value = 'Hortaliza'

# --- New row 265: Repollo, Crespo record, Segunda --------------------------
$ws.Range("A265").Value = 11
$ws.Range("B265").Value = 'Vega Monumental Concepción'
$ws.Range("C265").Value = 'Bíobío'
$ws.Range("D265").Value = 44736
$ws.Range("E265").Value = 8
$ws.Range("F265").Value = 100112006
$ws.Range("G265").Value = 'Repollo'
$ws.Range("H265").Value = 'Crespo record'
$ws.Range("I265").Value = 'Segunda'
$ws.Range("J265").Value = 1000
$ws.Range("K265").Value = 700
$ws.Range("L265").Value = 700
$ws.Range("M265").Value = 700
$ws.Range("N265").Value = '$/unidad'
$ws.Range("O265").Value = 'Región Metropolitana'
$ws.Range("P265").Value = 700
$ws.Range("Q265").Value = 1
$ws.Range("R265").Value = 'Hortaliza'
